$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-8 and 10-14
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 3
$ws.Range("F8").Value = -4
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = 11
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = 1
